$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at position 196 (shifts existing rows 196:339 down to 197:340)
$ws.Rows.Item(196).Insert()

# Populate the newly inserted row with the new record
$ws.Range("A196").Value = 4
$ws.Range("B196").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C196").Value = "Los Lagos"
$ws.Range("D196").Value = 44827
$ws.Range("E196").Value = 10
$ws.Range("F196").Value = 100112037
$ws.Range("G196").Value = "Cebollín"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 180
$ws.Range("K196").Value = 9000
$ws.Range("L196").Value = 9000
$ws.Range("M196").Value = 9000
$ws.Range("N196").Value = "`$/paquete 36 unidades"
$ws.Range("O196").Value = "Región Metropolitana"
$ws.Range("P196").Value = 250
$ws.Range("Q196").Value = 36
$ws.Range("R196").Value = "Hortaliza"
